$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells for the season record: Wins, Losses, Ties
# in the columns right after the existing last column (AB -> AC, AD, AE)
$ws.Range("AC1").Value = "Wins"
$ws.Range("AD1").Value = "Losses"
$ws.Range("AE1").Value = "Ties"

# Copy the header formatting (bold, centered, bordered) from the
# neighboring header cell so the new headers match the existing style
$ws.Range("AB1").Copy()
$ws.Range("AC1:AE1").PasteSpecial(-4122)

# Fill in the season record for every player row: every row shares the
# same team record (80 wins, 82 losses, 0 ties)
for ($r = 2; $r -le 48; $r++) {
    $ws.Cells.Item($r, 29).Value = 80
    $ws.Cells.Item($r, 30).Value = 82
    $ws.Cells.Item($r, 31).Value = 0
}
